# Restore the ShiftBookExcel2 header row on Sheet1 (A1:G1).
# Write B1:G1 before A1 so the shared-strings table is built in the
# same order ("studentID, name, date, shift, LIC, LIC vertified, timeStamp")
# as the target workbook, with A1's "timeStamp" landing last (index 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "studentID"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "shift"
$ws.Range("F1").Value = "LIC"
$ws.Range("G1").Value = "LIC vertified"
$ws.Range("A1").Value = "timeStamp"
